# xleash: AMIDST FIX element-dive loc-reporting.
#
# Adds a new "PevalAll" row to the first worksheet ("2"), inserted right
# after the existing "P-eval" row (row 9), pushing the remaining rows
# (the old row 11 "No Recurse" block and the table at rows 13-16) down
# by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row at row 10. Everything currently at row 10 and
# below (rows 11, 13-16) shifts down by one row (-> 11, 14-17).
$ws.Rows.Item(10).EntireRow.Insert() | Out-Null

# Populate the newly inserted row 10 with the PevalAll example.
$ws.Range("B10").Value = "PevalAll"
$ws.Range("C10").Value = '#eval sheet!::{"func": "pipe", "kwds":{"lax": false}, "args":[["df", {"index_col": null}], ["pyeval", {"include": "EVAL_COL", "eval_all": true}], "recurse"]}'

# Reflect the new selection (mirrors the author moving to the new row).
$ws.Range("C10").Select() | Out-Null
